$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text while writing, so numeric-looking
# strings like "1.000" or "0.3830" keep their exact text representation
# instead of being coerced into Double values by the COM Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.179.00"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "1.782.70"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "338.56"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.3830"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").Value = "0.3436"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "47.09"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "1.148"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "0.07376"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "23.16"
$ws.Range("E12").Value = "  +8.01%  "
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "6.456"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "7.380"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").Value = "1.785.33"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "0.00001074"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "0.06682"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "82.22"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "17.50"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "6.458"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "28.213.34"
$ws.Range("E23").Value = "  +3.94%  "
$ws.Range("D24").Value = "12.05"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").Value = "20.66"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").Value = "2.410"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").Value = "154.13"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "136.28"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "1.984.32"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "6.117"
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("D33").Value = "3.949"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").Value = "0.08881"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "0.02433"
$ws.Range("E36").Value = "  +4.14%  "
$ws.Range("D37").Value = "0.6847"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "5.328"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "0.06351"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").Value = "0.2172"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "1.246"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "1.496"
$ws.Range("E42").Value = "  -7.54%  "
$ws.Range("D43").Value = "8.265"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "14.12"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "0.6295"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "133.28"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "2.085"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").Value = "0.07441"
$ws.Range("E50").Value = "  +4.67%  "
$ws.Range("D51").Value = "1.209"
$ws.Range("E51").Value = "  +8.49%  "

# Restore the original (default/Normal) style on column D so the saved
# workbook keeps the same cell formatting as before (no explicit style
# index / number format applied), matching the source data layout.
$ws.Range("D2:D51").Style = "Normal"
